$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fix photo paths: replace ".jpg" with ".webp" for rows 2-30 (column C) ---
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    if ($old -ne $null -and $old.Contains(".jpg")) {
        $cell.Value = $old.Replace(".jpg", ".webp")
    }
}

# --- 2. Replace stale Google Drive links with the correct static webp paths ---
$ws.Cells.Item(31, 3).Value = "/static/images/profile_photos/005/VEC-005-02-199.webp"
$ws.Cells.Item(32, 3).Value = "/static/images/profile_photos/005/VEC-005-04-196.webp"

# --- 3. Adjust column widths (B, C, D, J) ---
$ws.Columns.Item(2).ColumnWidth = 18.333333333333336
$ws.Columns.Item(3).ColumnWidth = 52.33333333333333
$ws.Columns.Item(4).ColumnWidth = 15.0
$ws.Columns.Item(10).ColumnWidth = 17.0

# --- 4. Normalize row heights back to default (18.75) for header/data rows that previously had large custom heights ---
$rowsToReset = @(1,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
foreach ($r in $rowsToReset) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

# --- 5. Row 32 grows to fit the long new content ---
$ws.Rows.Item(32).RowHeight = 382.5
